$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item("Elektronarzedzia")
$ws2.Range("I1").Value = "ID"
$ws2.Range("J1").Value = "Nazwa"
$ws2.Range("K1").Value = "Opis"
$ws2.Range("L1").Value = "Typ ostrza"
$ws2.Range("M1").Value = "Moc silnika"
$ws2.Range("N1").Value = "Typ silnika"
$ws2.Range("O1").Value = "Typ zasilania"
$ws2.Range("P1").Value = "Cena"

$ws2.Range("I2").Value = "EPT165PP"
$ws2.Range("J2").Value = "PODSTAWOWA PILARKA TARCZOWA, PRZEWODOWA - 165 MM, 800 W"
$ws2.Range("K2").Value = "Podstawowa pilarka tarczowa do tarcz o rozmiarze 165 mm."
$ws2.Range("L2").Value = "PIŁA TARCZOWA"
$ws2.Range("M2").Value = "800 W"
$ws2.Range("N2").Value = "SZCZOTKOWY"
$ws2.Range("O2").Value = "PRZEWODOWY"
$ws2.Range("P2").NumberFormat = "@"
$ws2.Range("P2").Value = "300.0"

$ws2.Range("I3").Value = "EPT165UP"
$ws2.Range("J3").Value = "UNIWERSALNA PILARKA TARCZOWA, PRZEWODOWA - 165 MM, 1200 W"
$ws2.Range("K3").Value = "Uniwersalna pilarka tarczowa do tarcz o rozmiarze 165 mm."
$ws2.Range("L3").Value = "PIŁA TARCZOWA"
$ws2.Range("M3").Value = "1200 W"
$ws2.Range("N3").Value = "SZCZOTKOWY"
$ws2.Range("O3").Value = "PRZEWODOWY"
$ws2.Range("P3").NumberFormat = "@"
$ws2.Range("P3").Value = "500.0"

$ws2.Range("I4").Value = "EPT165ZP"
$ws2.Range("J4").Value = "ZAAWANSOWANA PILARKA TARCZOWA, PRZEWODOWA - 165 MM, 1600 W"
$ws2.Range("K4").Value = "Zaawansowana pilarka tarczowa do tarcz o rozmiarze 165 mm."
$ws2.Range("L4").Value = "PIŁA TARCZOWA"
$ws2.Range("M4").Value = "1600 W"
$ws2.Range("N4").Value = "BEZSZCZOTKOWY"
$ws2.Range("O4").Value = "PRZEWODOWY"
$ws2.Range("P4").NumberFormat = "@"
$ws2.Range("P4").Value = "800.0"

$ws2.Range("I5").Value = "EPT165PA"
$ws2.Range("J5").Value = "PODSTAWOWA PILARKA TARCZOWA, AKUMULATOROWA - 165 MM, 800 W"
$ws2.Range("K5").Value = "Podstawowa, akumulatorowa pilarka tarczowa do tarcz o rozmiarze 165 mm."
$ws2.Range("L5").Value = "PIŁA TARCZOWA"
$ws2.Range("M5").Value = "800 W"
$ws2.Range("N5").Value = "SZCZOTKOWY"
$ws2.Range("O5").Value = "AKUMULATOR 18V"
$ws2.Range("P5").NumberFormat = "@"
$ws2.Range("P5").Value = "420.0"

$ws2.Range("I6").Value = "EPT165UA"
$ws2.Range("J6").Value = "UNIWERSALNA PILARKA TARCZOWA, AKUMULATOROWA - 165 MM, 1200 W"
$ws2.Range("K6").Value = "Uniwersalna, akumulatorowa pilarka tarczowa do tarcz o rozmiarze 165 mm."
$ws2.Range("L6").Value = "PIŁA TARCZOWA"
$ws2.Range("M6").Value = "1200 W"
$ws2.Range("N6").Value = "SZCZOTKOWY"
$ws2.Range("O6").Value = "AKUMULATOR 18V"
$ws2.Range("P6").NumberFormat = "@"
$ws2.Range("P6").Value = "650.0"

$ws2.Range("I7").Value = "EPT165ZA"
$ws2.Range("J7").Value = "ZAAWANSOWANA PILARKA TARCZOWA, AKUMULATOROWA - 165 MM, 1600 W"
$ws2.Range("K7").Value = "Zaawansowana, akumulatorowa pilarka tarczowa do tarcz o rozmiarze 165 mm."
$ws2.Range("L7").Value = "PIŁA TARCZOWA"
$ws2.Range("M7").Value = "1600 W"
$ws2.Range("N7").Value = "BEZSZCZOTKOWY"
$ws2.Range("O7").Value = "AKUMULATOR 18V"
$ws2.Range("P7").NumberFormat = "@"
$ws2.Range("P7").Value = "990.0"

$ws2.Range("I8").Value = "EPSUP"
$ws2.Range("J8").Value = "UNIWERSALNA PIŁA SZABLASTA, 800 W"
$ws2.Range("K8").Value = "Uniwersalna piła szablasta o mocy 800 W"
$ws2.Range("L8").Value = "BRZESZCZOT BAGNETOWY"
$ws2.Range("M8").Value = "800 W"
$ws2.Range("N8").Value = "SZCZOTKOWY"
$ws2.Range("O8").Value = "PRZEWODOWY"
$ws2.Range("P8").NumberFormat = "@"
$ws2.Range("P8").Value = "320.0"

$ws2.Range("I9").Value = "EPSZP"
$ws2.Range("J9").Value = "ZAAWANSOWANA PIŁA SZABLASTA, 1200 W"
$ws2.Range("K9").Value = "Zaawansowana piła szablasta o mocy 1200 W"
$ws2.Range("L9").Value = "BRZESZCZOT BAGNETOWY"
$ws2.Range("M9").Value = "1200 W"
$ws2.Range("N9").Value = "BEZSZCZOTKOWY"
$ws2.Range("O9").Value = "PRZEWODOWY"
$ws2.Range("P9").NumberFormat = "@"
$ws2.Range("P9").Value = "690.0"

$ws2.Range("I10").Value = "EPSUA"
$ws2.Range("J10").Value = "UNIWERSALNA PIŁA SZABLASTA, AKUMULATOROWA , 800 W"
$ws2.Range("K10").Value = "Uniwersalna piła szablasta o mocy 800 W, akumulatorowa"
$ws2.Range("L10").Value = "BRZESZCZOT BAGNETOWY"
$ws2.Range("M10").Value = "800 W"
$ws2.Range("N10").Value = "SZCZOTKOWY"
$ws2.Range("O10").Value = "AKUMULATOR 18V"
$ws2.Range("P10").NumberFormat = "@"
$ws2.Range("P10").Value = "520.0"

$ws2.Range("I11").Value = "EPSZA"
$ws2.Range("J11").Value = "ZAAWANSOWANA PIŁA SZABLASTA, AKUMULATOROWA , 1200 W"
$ws2.Range("K11").Value = "Zaawansowana piła szablasta o mocy 1200 W, akumulatorowa"
$ws2.Range("L11").Value = "BRZESZCZOT BAGNETOWY"
$ws2.Range("M11").Value = "1200 W"
$ws2.Range("N11").Value = "BEZSZCZOTKOWY"
$ws2.Range("O11").Value = "AKUMULATOR 18V"
$ws2.Range("P11").NumberFormat = "@"
$ws2.Range("P11").Value = "870.0"

$ws3 = $wb.Worksheets.Item("Ostrza")
$ws3.Range("L1").Value = "ID"
$ws3.Range("M1").Value = "Nazwa"
$ws3.Range("N1").Value = "Typ"
$ws3.Range("O1").Value = "Opis"
$ws3.Range("P1").Value = "Dlugosc"
$ws3.Range("Q1").Value = "Srednica"
$ws3.Range("R1").Value = "Grubosc"
$ws3.Range("S1").Value = "Material"
$ws3.Range("T1").Value = "Liczba zebow (lub na cal)"
$ws3.Range("U1").Value = "Zastosowanie"
$ws3.Range("V1").Value = "Cena"

$ws3.Range("L2").Value = "PTDP165"
$ws3.Range("M2").Value = "PIŁA TARCZOWA DO DREWNA - PODSTAWOWA 165"
$ws3.Range("N2").Value = "PIŁA TARCZOWA"
$ws3.Range("O2").Value = "Piła tarczowa o rozmiarze 165 mm do podstawowych zastowań przy drewnie."
$ws3.Range("P2").NumberFormat = "@"
$ws3.Range("P2").Value = "None"
$ws3.Range("Q2").NumberFormat = "@"
$ws3.Range("Q2").Value = "165.0"
$ws3.Range("R2").NumberFormat = "@"
$ws3.Range("R2").Value = "1.5"
$ws3.Range("S2").Value = "WĘGLIK SPIEKANY"
$ws3.Range("T2").NumberFormat = "@"
$ws3.Range("T2").Value = "24.0"
$ws3.Range("U2").Value = "DREWNO`nFORNIR`nPŁYTY WIÓROWE"
$ws3.Range("V2").NumberFormat = "@"
$ws3.Range("V2").Value = "100.0"

$ws3.Range("L3").Value = "PTDZ165"
$ws3.Range("M3").Value = "PIŁA TARCZOWA DO DREWNA - ZAAWANSOWANA 165"
$ws3.Range("N3").Value = "PIŁA TARCZOWA"
$ws3.Range("O3").Value = "Piła tarczowa do drewna o rozmiarze 165 mm  do najtrudniejszych zadań."
$ws3.Range("P3").NumberFormat = "@"
$ws3.Range("P3").Value = "None"
$ws3.Range("Q3").NumberFormat = "@"
$ws3.Range("Q3").Value = "165.0"
$ws3.Range("R3").NumberFormat = "@"
$ws3.Range("R3").Value = "1.5"
$ws3.Range("S3").Value = "WĘGLIK SPIEKANY"
$ws3.Range("T3").NumberFormat = "@"
$ws3.Range("T3").Value = "24.0"
$ws3.Range("U3").Value = "DREWNO`nFORNIR`nPŁYTY WIÓROWE"
$ws3.Range("V3").NumberFormat = "@"
$ws3.Range("V3").Value = "180.0"

$ws3.Range("L4").Value = "PTDP210"
$ws3.Range("M4").Value = "PIŁA TARCZOWA DO DREWNA -  PODSTAWOWA 210"
$ws3.Range("N4").Value = "PIŁA TARCZOWA"
$ws3.Range("O4").Value = "Piła tarczowa o rozmiarze 210 mm  do podstawowych zastowań przy drewnie."
$ws3.Range("P4").NumberFormat = "@"
$ws3.Range("P4").Value = "None"
$ws3.Range("Q4").NumberFormat = "@"
$ws3.Range("Q4").Value = "210.0"
$ws3.Range("R4").NumberFormat = "@"
$ws3.Range("R4").Value = "1.9"
$ws3.Range("S4").Value = "WĘGLIK SPIEKANY"
$ws3.Range("T4").NumberFormat = "@"
$ws3.Range("T4").Value = "24.0"
$ws3.Range("U4").Value = "DREWNO`nFORNIR`nPŁYTY WIÓROWE"
$ws3.Range("V4").NumberFormat = "@"
$ws3.Range("V4").Value = "140.0"

$ws3.Range("L5").Value = "PTDZ210"
$ws3.Range("M5").Value = "PIŁA TARCZOWA DO DREWNA - ZAAWANSOWANA 210"
$ws3.Range("N5").Value = "PIŁA TARCZOWA"
$ws3.Range("O5").Value = "Piła tarczowa do drewna  o rozmiarze 210 mm do najtrudniejszych zadań."
$ws3.Range("P5").NumberFormat = "@"
$ws3.Range("P5").Value = "None"
$ws3.Range("Q5").NumberFormat = "@"
$ws3.Range("Q5").Value = "210.0"
$ws3.Range("R5").NumberFormat = "@"
$ws3.Range("R5").Value = "1.9"
$ws3.Range("S5").Value = "WĘGLIK SPIEKANY"
$ws3.Range("T5").NumberFormat = "@"
$ws3.Range("T5").Value = "24.0"
$ws3.Range("U5").Value = "DREWNO`nFORNIR`nPŁYTY WIÓROWE"
$ws3.Range("V5").NumberFormat = "@"
$ws3.Range("V5").Value = "220.0"

$ws3.Range("L6").Value = "PTU165"
$ws3.Range("M6").Value = "PIŁA TARCZOWA - UNIWERSALNA 165"
$ws3.Range("N6").Value = "PIŁA TARCZOWA"
$ws3.Range("O6").Value = "Uniwersalna piła tarczowa o rozmiarze 165 mm do roznych materiałów."
$ws3.Range("P6").NumberFormat = "@"
$ws3.Range("P6").Value = "None"
$ws3.Range("Q6").NumberFormat = "@"
$ws3.Range("Q6").Value = "165.0"
$ws3.Range("R6").NumberFormat = "@"
$ws3.Range("R6").Value = "1.5"
$ws3.Range("S6").Value = "WĘGLIK SPIEKANY"
$ws3.Range("T6").NumberFormat = "@"
$ws3.Range("T6").Value = "48.0"
$ws3.Range("U6").Value = "ALUMINIUM`nDREWNO`nLAMINAT`nPVC"
$ws3.Range("V6").NumberFormat = "@"
$ws3.Range("V6").Value = "210.0"

$ws3.Range("L7").Value = "PTU210"
$ws3.Range("M7").Value = "PIŁA TARCZOWA - UNIWERSALNA 210"
$ws3.Range("N7").Value = "PIŁA TARCZOWA"
$ws3.Range("O7").Value = "Uniwersalna piła tarczowa o rozmiarze 165 mm do roznych materiałów."
$ws3.Range("P7").NumberFormat = "@"
$ws3.Range("P7").Value = "None"
$ws3.Range("Q7").NumberFormat = "@"
$ws3.Range("Q7").Value = "210.0"
$ws3.Range("R7").NumberFormat = "@"
$ws3.Range("R7").Value = "1.5"
$ws3.Range("S7").Value = "WĘGLIK SPIEKANY"
$ws3.Range("T7").NumberFormat = "@"
$ws3.Range("T7").Value = "48.0"
$ws3.Range("U7").Value = "ALUMINIUM`nDREWNO`nLAMINAT`nPVC"
$ws3.Range("V7").NumberFormat = "@"
$ws3.Range("V7").Value = "270.0"

$ws3.Range("L8").Value = "BBM150"
$ws3.Range("M8").Value = "BRZESZCZOT BAGNETOWY - DO METALU 150"
$ws3.Range("N8").Value = "BRZESZCZOT BAGNETOWY"
$ws3.Range("O8").Value = "Brzeszczot bagnetowy do metalu o długości 150 mm."
$ws3.Range("P8").NumberFormat = "@"
$ws3.Range("P8").Value = "150.0"
$ws3.Range("Q8").NumberFormat = "@"
$ws3.Range("Q8").Value = "None"
$ws3.Range("R8").NumberFormat = "@"
$ws3.Range("R8").Value = "0.9"
$ws3.Range("S8").Value = "HSS"
$ws3.Range("T8").NumberFormat = "@"
$ws3.Range("T8").Value = "24.0"
$ws3.Range("U8").Value = "LAMINAT`nMETAL`nPVC"
$ws3.Range("V8").NumberFormat = "@"
$ws3.Range("V8").Value = "15.0"

$ws3.Range("L9").Value = "BBM230"
$ws3.Range("M9").Value = "BRZESZCZOT BAGNETOWY - DO METALU 230"
$ws3.Range("N9").Value = "BRZESZCZOT BAGNETOWY"
$ws3.Range("O9").Value = "Brzeszczot bagnetowy do metalu o długości 230 mm."
$ws3.Range("P9").NumberFormat = "@"
$ws3.Range("P9").Value = "230.0"
$ws3.Range("Q9").NumberFormat = "@"
$ws3.Range("Q9").Value = "None"
$ws3.Range("R9").NumberFormat = "@"
$ws3.Range("R9").Value = "0.9"
$ws3.Range("S9").Value = "HSS"
$ws3.Range("T9").NumberFormat = "@"
$ws3.Range("T9").Value = "24.0"
$ws3.Range("U9").Value = "LAMINAT`nMETAL`nPVC"
$ws3.Range("V9").NumberFormat = "@"
$ws3.Range("V9").Value = "25.0"

$ws3.Range("L10").Value = "BBD150"
$ws3.Range("M10").Value = "BRZESZCZOT BAGNETOWY - DO DREWNA 150"
$ws3.Range("N10").Value = "BRZESZCZOT BAGNETOWY"
$ws3.Range("O10").Value = "Brzeszczot bagnetowy do drewna o długości 150 mm."
$ws3.Range("P10").NumberFormat = "@"
$ws3.Range("P10").Value = "150.0"
$ws3.Range("Q10").NumberFormat = "@"
$ws3.Range("Q10").Value = "None"
$ws3.Range("R10").NumberFormat = "@"
$ws3.Range("R10").Value = "1.3"
$ws3.Range("S10").Value = "BIMETAL"
$ws3.Range("T10").NumberFormat = "@"
$ws3.Range("T10").Value = "8.0"
$ws3.Range("U10").Value = "DREWNO`nFORNIR`nPŁYTY WIÓROWE"
$ws3.Range("V10").NumberFormat = "@"
$ws3.Range("V10").Value = "15.0"

$ws3.Range("L11").Value = "BBD230"
$ws3.Range("M11").Value = "BRZESZCZOT BAGNETOWY - DO DREWNA 230"
$ws3.Range("N11").Value = "BRZESZCZOT BAGNETOWY"
$ws3.Range("O11").Value = "Brzeszczot bagnetowy do drewna o długości 230 mm."
$ws3.Range("P11").NumberFormat = "@"
$ws3.Range("P11").Value = "230.0"
$ws3.Range("Q11").NumberFormat = "@"
$ws3.Range("Q11").Value = "None"
$ws3.Range("R11").NumberFormat = "@"
$ws3.Range("R11").Value = "1.3"
$ws3.Range("S11").Value = "BIMETAL"
$ws3.Range("T11").NumberFormat = "@"
$ws3.Range("T11").Value = "8.0"
$ws3.Range("U11").Value = "DREWNO`nFORNIR`nPŁYTY WIÓROWE"
$ws3.Range("V11").NumberFormat = "@"
$ws3.Range("V11").Value = "25.0"

$ws3.Range("L12").Value = "BBM300"
$ws3.Range("M12").Value = "BRZESZCZOT BAGNETOWY - DO METALU 300"
$ws3.Range("N12").Value = "BRZESZCZOT BAGNETOWY"
$ws3.Range("O12").Value = "Brzeszczot bagnetowy do metalu o długości 300 mm."
$ws3.Range("P12").NumberFormat = "@"
$ws3.Range("P12").Value = "230.0"
$ws3.Range("Q12").NumberFormat = "@"
$ws3.Range("Q12").Value = "None"
$ws3.Range("R12").NumberFormat = "@"
$ws3.Range("R12").Value = "0.9"
$ws3.Range("S12").Value = "HSS"
$ws3.Range("T12").NumberFormat = "@"
$ws3.Range("T12").Value = "24.0"
$ws3.Range("U12").Value = "LAMINAT`nMETAL`nPVC"
$ws3.Range("V12").NumberFormat = "@"
$ws3.Range("V12").Value = "30.0"

$ws3.Range("L13").Value = "BBD300"
$ws3.Range("M13").Value = "BRZESZCZOT BAGNETOWY - DO DREWNA 300"
$ws3.Range("N13").Value = "BRZESZCZOT BAGNETOWY"
$ws3.Range("O13").Value = "Brzeszczot bagnetowy do drewna o długości 300 mm."
$ws3.Range("P13").NumberFormat = "@"
$ws3.Range("P13").Value = "230.0"
$ws3.Range("Q13").NumberFormat = "@"
$ws3.Range("Q13").Value = "None"
$ws3.Range("R13").NumberFormat = "@"
$ws3.Range("R13").Value = "1.3"
$ws3.Range("S13").Value = "BIMETAL"
$ws3.Range("T13").NumberFormat = "@"
$ws3.Range("T13").Value = "8.0"
$ws3.Range("U13").Value = "DREWNO`nFORNIR`nPŁYTY WIÓROWE"
$ws3.Range("V13").NumberFormat = "@"
$ws3.Range("V13").Value = "30.0"
